# Fill in the previously-blank "Mid Paper" marks (columns D-G) for three
# students on the "Senior Six" sheet that had been left empty by mistake.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Senior Six")

# Row 7 - OCHEN ATIKU HUSSEIN
$ws.Range("D7:G7").WrapText = $false
$ws.Range("D7").Value = 75.0
$ws.Range("E7").Value = 70.0
$ws.Range("F7").Value = 86.0
$ws.Range("G7").Value = 66.0

# Row 11 - OKWERA ERICK
$ws.Range("D11:G11").WrapText = $false
$ws.Range("D11").Value = 68.0
$ws.Range("E11").Value = 54.0
$ws.Range("F11").Value = 65.0
$ws.Range("G11").Value = 54.0

# Row 14 - OPIYO ATIKU HASSAN
$ws.Range("D14:G14").WrapText = $false
$ws.Range("D14").Value = 72.0
$ws.Range("E14").Value = 68.0
$ws.Range("F14").Value = 70.0
$ws.Range("G14").Value = 54.0
